# Add two new locator rows to the "locators" sheet:
#   row 14 - "all courses link" header locator
#   row 15 - "selenium webdriver 4 with java" course link locator
# and move the active selection to E15 (the last filled cell), matching
# the author's workflow of filling the new rows then leaving the
# selection on the trailing status cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the same "Good" (green) style already applied to the other
# Status cells in column E.
$statusStyle = $ws.Range("E13").Style

# --- Row 14: header / all courses link -------------------------------
$ws.Cells.Item(14, 2).Value = "all courses link"
$ws.Cells.Item(14, 1).Value = "header"
$ws.Cells.Item(14, 3).Value = "xpath"
$ws.Cells.Item(14, 4).Value = "//a[contains(@href,'/courses') and contains(text(),'ALL COURSES')]"
$ws.Cells.Item(14, 5).Style = $statusStyle

# --- Row 15: all courses / selenium webdriver 4 with java -------------
$ws.Cells.Item(15, 1).Value = "all courses"
$ws.Cells.Item(15, 2).Value = "selenium webdriver 4 with java"
$ws.Cells.Item(15, 3).Value = "xpath"
$ws.Cells.Item(15, 4).Value = "//a[@href='/courses/selenium-webdriver-with-java']"
$ws.Cells.Item(15, 5).Style = $statusStyle

# Update the selection to match the author's final cursor position.
$ws.Range("E15").Select()
